# "cambios de agosto, puntos fe de ratas e historico"
#
# Update the reporting-period figures on the visible sheet ("Reporte de
# Formatos") row 8: the report moves from the 2021 Q3/Q4 window to the
# 2022 Q1/Q2 window, and the "last updated" dates move forward too.
# Also update the window's scroll position / active selection to match
# where the author left the cursor (column U in view, W10 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# --- Row 8 data changes -----------------------------------------------
$ws.Range("A8").Value = 2022        # Ejercicio: 2021 -> 2022
$ws.Range("B8").Value = 44562       # Fecha de inicio del periodo que se informa
$ws.Range("C8").Value = 44742       # Fecha de término del periodo que se informa
$ws.Range("U8").Value = 44753       # Fecha de validación
$ws.Range("V8").Value = 44753       # Fecha de actualización

# --- View / selection changes ------------------------------------------
$ws.Activate()

# Best effort: scroll the window so column U is the left-most visible
# column (mirrors topLeftCell="U2" in the saved sheetView).
try {
    $excel.ActiveWindow.ScrollColumn = $ws.Range("U2").Column
    $excel.ActiveWindow.ScrollRow = $ws.Range("U2").Row
} catch {
}

# Select W10, matching the author's final selection/cursor position.
$ws.Range("W10").Select()
